# Refresh the cryptos price table (Coinranking snapshot) in place.
# Mirrors the automated "Updated cryptos list" GitHub Actions commit:
# new Price/Volume(1h) figures for every coin, plus a couple of rank
# swaps where two coins traded places (rows 33/34 and 43-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.833.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "'3.109.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'525.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "'141.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.107.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "'0.433"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").Value = "'0.384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("D13").Value = "'3.642.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "'26.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "'57.883.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "'3.107.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'12.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'8.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "'337.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'0.508"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "'66.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "'0.168"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'0.0₃0927"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "'6.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.00%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").Value = "'1.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'20.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'153.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'4.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("D37").Value = "'6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").Value = "'26.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'1.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").Value = "'0.0667"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'3.148.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").Value = "'0.685"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'36.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'2.296.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "'0.0258"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "'0.989"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.32%  "
$ws.Range("D50").Value = "'20.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").Value = "'5.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.54%  "
